$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings are not
# reinterpreted as numbers (matches the original inlineStr/text cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.177.72'
$ws.Range("E2").Value = '  +1.67%  '

$ws.Range("D3").Value = '3.316.98'
$ws.Range("E3").Value = '  +6.14%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '601.29'
$ws.Range("E5").Value = '  +1.14%  '

$ws.Range("D6").Value = '144.01'
$ws.Range("E6").Value = '  +5.55%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '3.312.49'
$ws.Range("E8").Value = '  +6.20%  '

$ws.Range("E9").Value = '  +0.60%  '

$ws.Range("E10").Value = '  +2.30%  '

$ws.Range("D11").Value = '5.48'
$ws.Range("E11").Value = '  +2.67%  '

$ws.Range("E12").Value = '  +3.31%  '

$ws.Range("E13").Value = '  +0.14%  '

$ws.Range("D14").Value = '35.05'
$ws.Range("E14").Value = '  +2.64%  '

$ws.Range("D15").Value = '3.869.78'
$ws.Range("E15").Value = '  +6.33%  '

$ws.Range("D16").Value = '0.120'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").Value = '3.317.62'
$ws.Range("E17").Value = '  +6.52%  '

$ws.Range("D18").Value = '64.264.47'
$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").Value = '6.92'
$ws.Range("E19").Value = '  +3.00%  '

$ws.Range("D20").Value = '485.46'
$ws.Range("E20").Value = '  +1.74%  '

$ws.Range("D21").Value = '14.38'
$ws.Range("E21").Value = '  +1.35%  '

$ws.Range("D22").Value = '0.746'
$ws.Range("E22").Value = '  +6.67%  '

$ws.Range("E23").Value = '  +5.61%  '

$ws.Range("D24").Value = '13.62'
$ws.Range("E24").Value = '  +3.77%  '

$ws.Range("D25").Value = '84.96'
$ws.Range("E25").Value = '  -3.24%  '

$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  +2.78%  '

$ws.Range("D28").Value = '8.33'
$ws.Range("E28").Value = '  +4.04%  '

$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").Value = '7.28'
$ws.Range("E29").Value = '  +0.82%  '

$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = '2.17'
$ws.Range("E31").Value = '  +4.70%  '

$ws.Range("D32").Value = '28.48'
$ws.Range("E32").Value = '  +4.39%  '

$ws.Range("D33").Value = '0.108'
$ws.Range("E33").Value = '  -1.43%  '

$ws.Range("E34").Value = '  +1.60%  '

$ws.Range("E35").Value = '  +2.48%  '

$ws.Range("D36").Value = '6.04'
$ws.Range("E36").Value = '  +3.05%  '

$ws.Range("D37").Value = '53.39'
$ws.Range("E37").Value = '  +2.56%  '

$ws.Range("D38").Value = '0.0₃0741'
$ws.Range("E38").Value = '  +3.78%  '

$ws.Range("E39").Value = '  +3.07%  '

$ws.Range("D40").Value = '432.88'
$ws.Range("E40").Value = '  +2.69%  '

$ws.Range("D41").Value = '2.81'
$ws.Range("E41").Value = '  +3.60%  '

$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").Value = '8.50'
$ws.Range("E42").Value = '  +2.81%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.020.81'
$ws.Range("E43").Value = '  +5.14%  '

$ws.Range("D44").Value = '0.111'
$ws.Range("E44").Value = '  -5.17%  '

$ws.Range("E45").Value = '  +4.91%  '

$ws.Range("D46").Value = '2.25'
$ws.Range("E46").Value = '  +6.14%  '

$ws.Range("D47").Value = '26.44'
$ws.Range("E47").Value = '  +3.31%  '

$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("E49").Value = '  +2.52%  '

$ws.Range("E50").Value = '  +1.76%  '

$ws.Range("D51").Value = '35.46'
$ws.Range("E51").Value = '  +15.22%  '
